$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("110:110").Insert()

$ws.Range("A110").Value = 4
$ws.Range("B110").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C110").Value = 'Los Lagos'
$ws.Range("D110").Value = 44603
$ws.Range("D110").NumberFormat = $ws.Range("D111").NumberFormat
$ws.Range("E110").Value = 10
$ws.Range("F110").Value = 100112039
$ws.Range("G110").Value = 'Ciboulette'
$ws.Range("H110").Value = 'Sin especificar'
$ws.Range("I110").Value = 'Primera'
$ws.Range("J110").Value = 240
$ws.Range("K110").Value = 2500
$ws.Range("L110").Value = 2500
$ws.Range("M110").Value = 2500
$ws.Range("N110").Value = '$/docena de atados'
$ws.Range("O110").Value = 'Región Metropolitana'
$ws.Range("P110").Value = 833
$ws.Range("Q110").Value = 3
$ws.Range("R110").Value = 'Hortaliza'
